$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 3645.3684
$ws.Range("I62").Value = 3687.5557
$ws.Range("J62").Value = 3607.4
$ws.Range("K62").Value = 3687.5557
$ws.Range("L62").Value = 3607.4
$ws.Range("M62").Value = -3063.5557
$ws.Range("N62").Value = -4855.4
$ws.Range("H65").Value = 3645.3684
$ws.Range("I65").Value = 3687.5557
$ws.Range("J65").Value = 3607.4
$ws.Range("K65").Value = 18437.7785
$ws.Range("L65").Value = 18037
$ws.Range("M65").Value = -15317.7785
$ws.Range("N65").Value = -24277
$ws.Range("H111").Value = 2877.3157
$ws.Range("I111").Value = 4667.3335
$ws.Range("J111").Value = 1266.3
$ws.Range("K111").Value = 14002.0005
$ws.Range("L111").Value = 3798.9
$ws.Range("M111").Value = -10935.0005
$ws.Range("N111").Value = -9932.9
$ws.Range("H116").Value = 2015.317
$ws.Range("I116").Value = 1974.9143
$ws.Range("J116").Value = 2251
$ws.Range("K116").Value = 1974.9143
$ws.Range("L116").Value = 2251
$ws.Range("M116").Value = 1467.0857
$ws.Range("N116").Value = -9135
$ws.Range("H137").Value = 26361.854
$ws.Range("I137").Value = 1551.1666
$ws.Range("J137").Value = 94027.37
$ws.Range("K137").Value = 4653.4998
$ws.Range("L137").Value = 282082.11
$ws.Range("M137").Value = -2103.4998
$ws.Range("N137").Value = -287182.11

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 19236208
$ws.Range("I32").Value = 21278592
$ws.Range("J32").Value = 37800
$ws.Range("K32").Value = 21278592
$ws.Range("L32").Value = 37800
$ws.Range("M32").Value = -21278305
$ws.Range("N32").Value = -38374
$ws.Range("H74").Value = 4836523.5
$ws.Range("I74").Value = 7420.3125
$ws.Range("J74").Value = 15874473
$ws.Range("K74").Value = 7420.3125
$ws.Range("L74").Value = 15874473
$ws.Range("M74").Value = -6546.3125
$ws.Range("N74").Value = -15876221
$ws.Range("H77").Value = 4836523.5
$ws.Range("I77").Value = 7420.3125
$ws.Range("J77").Value = 15874473
$ws.Range("K77").Value = 37101.5625
$ws.Range("L77").Value = 79372365
$ws.Range("M77").Value = -32733.5625
$ws.Range("N77").Value = -79381101
$ws.Range("H110").Value = 1458.8
$ws.Range("I110").Value = 1476.16
$ws.Range("J110").Value = 1372
$ws.Range("K110").Value = 1476.16
$ws.Range("L110").Value = 1372
$ws.Range("M110").Value = 568.8399999999999
$ws.Range("N110").Value = -5462
$ws.Range("H122").Value = 1561.2307
$ws.Range("I122").Value = 1304.381
$ws.Range("J122").Value = 2640
$ws.Range("K122").Value = 3913.143
$ws.Range("L122").Value = 7920
$ws.Range("M122").Value = -1463.143
$ws.Range("N122").Value = -12820

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 1310.0209
$ws.Range("I134").Value = 1336.5814
$ws.Range("J134").Value = 1081.6
$ws.Range("K134").Value = 4009.7442
$ws.Range("L134").Value = 3244.8
$ws.Range("M134").Value = -1474.7442
$ws.Range("N134").Value = -8314.799999999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 2425
$ws.Range("I16").Value = 1898.5714
$ws.Range("J16").Value = 3162
$ws.Range("K16").Value = 1898.5714
$ws.Range("L16").Value = 3162
$ws.Range("M16").Value = -1611.5714
$ws.Range("N16").Value = -3736
$ws.Range("H113").Value = 2425
$ws.Range("I113").Value = 1898.5714
$ws.Range("J113").Value = 3162
$ws.Range("K113").Value = 1898.5714
$ws.Range("L113").Value = 3162
$ws.Range("M113").Value = 271.4286
$ws.Range("N113").Value = -7502
$ws.Range("H132").Value = 1948.8948
$ws.Range("I132").Value = 1732
$ws.Range("J132").Value = 2418.8333
$ws.Range("K132").Value = 5196
$ws.Range("L132").Value = 7256.499899999999
$ws.Range("M132").Value = -2666
$ws.Range("N132").Value = -12316.4999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 897.2857
$ws.Range("I122").Value = 731
$ws.Range("J122").Value = 1119
$ws.Range("K122").Value = 6579
$ws.Range("L122").Value = 10071
$ws.Range("M122").Value = -4129
$ws.Range("N122").Value = -14971

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 281.76923
$ws.Range("I107").Value = 216.4
$ws.Range("J107").Value = 499.66666
$ws.Range("K107").Value = 216.4
$ws.Range("L107").Value = 499.66666
$ws.Range("M107").Value = 1703.6
$ws.Range("N107").Value = -4339.66666
$ws.Range("H122").Value = 2859118
$ws.Range("I122").Value = 4547141
$ws.Range("J122").Value = 2463.077
$ws.Range("K122").Value = 13641423
$ws.Range("L122").Value = 7389.231000000001
$ws.Range("M122").Value = -13638973
$ws.Range("N122").Value = -12289.231

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 6974.4116
$ws.Range("I61").Value = 6974.4116
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 6974.4116
$ws.Range("L61").Value = 0
$ws.Range("M61").Value = -6772.4116
$ws.Range("N61").ClearContents()
$ws.Range("H113").Value = 6974.4116
$ws.Range("I113").Value = 6974.4116
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 6974.4116
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = -4804.4116
$ws.Range("N113").ClearContents()
$ws.Range("H122").Value = 1673.9375
$ws.Range("I122").Value = 1681.0834
$ws.Range("J122").Value = 1652.5
$ws.Range("K122").Value = 5043.2502
$ws.Range("L122").Value = 4957.5
$ws.Range("M122").Value = -2593.2502
$ws.Range("N122").Value = -9857.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 2304
$ws.Range("I96").Value = 1600
$ws.Range("J96").Value = 2480
$ws.Range("K96").Value = 1600
$ws.Range("L96").Value = 2480
$ws.Range("M96").Value = -227
$ws.Range("N96").Value = -5226
$ws.Range("H107").Value = 556.64
$ws.Range("I107").Value = 506.21738
$ws.Range("J107").Value = 1136.5
$ws.Range("K107").Value = 1518.65214
$ws.Range("L107").Value = 3409.5
$ws.Range("M107").Value = 401.3478600000001
$ws.Range("N107").Value = -7249.5
$ws.Range("H122").Value = 2505.6667
$ws.Range("I122").Value = 1979.1111
$ws.Range("J122").Value = 2821.6
$ws.Range("K122").Value = 5937.3333
$ws.Range("L122").Value = 8464.799999999999
$ws.Range("M122").Value = -3487.3333
$ws.Range("N122").Value = -13364.8
